$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap plasmid "pAGK050_R5" (50) and "pAGK103_R1" (103) data between row 3 and row 4
$ws.Range("B3").Value = "pAGK103_R1"
$ws.Range("C3").Value = 103
$ws.Range("F3").Value = "GTAGGACACA"
$ws.Range("G3").Value = "GGTTATCTGG"

$ws.Range("B4").Value = "pAGK050_R5"
$ws.Range("C4").Value = 50
$ws.Range("F4").Value = "GGGTAAAGGC"
$ws.Range("G4").Value = "GTACGGCATC"

# Clean up row 21: match formatting used by the data rows above and drop the
# stray F21/G21 cells entirely
$ws.Range("B21").Font.Name = "Calibri"
$ws.Range("B21").Font.Size = 11
$ws.Range("F21").Clear()
$ws.Range("G21").Clear()

# Move the active cell selection to C21
$ws.Range("C21").Select()
